$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (C and D) for "lat" and "lng", pushing the
# existing "type" and "bullets" columns from C/D to E/F.
$ws.Columns("C:D").Insert()

# Headers
$ws.Range("C1").Value = "lat"
$ws.Range("D1").Value = "lng"

# Fix a typo in the existing "type" column (now column E): Armado -> Armando
$ws.Range("E3").Value = "Armando"

# Populate the new lat/lng columns with coordinates for each city
# Row 2: Bordeaux, France
$ws.Range("C2").Value = 44.8378
$ws.Range("D2").Value = -0.5792

# Row 3: Leipzig, Germany
$ws.Range("C3").Value = 51.3397
$ws.Range("D3").Value = 12.3731

# Row 4: Fernando de Noronha, Brazil
$ws.Range("C4").Value = -3.857
$ws.Range("D4").Value = -32.429

# Remove the trailing empty, pre-formatted rows (5-13) that are no longer
# part of the data range.
$ws.Rows("5:13").Delete()
